$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (plain "In"/"Out" variable headers): x -> k(x/j), y -> j(y/k)
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "k(x/j)"
$ws.Range("C2").Value = "j(y/k)"
$ws.Range("D2").Value = "k(x/j)"
$ws.Range("E2").Value = "j(y/k)"
$ws.Range("F2").Value = "k(x/j)"
$ws.Range("G2").Value = "j(y/k)"

# ---------------------------------------------------------------------------
# Row 3: Sine
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "k+k*sin((pi/2)(x/j)-(pi/2))"
$ws.Range("C3").Value = "j + j(2*asin((y/k)-1))/pi"
$ws.Range("D3").Value = "k*sin((pi/2)(x/j))"
$ws.Range("E3").Value = "j(2*asin((y/k)))/pi"
$ws.Range("F3").Value = "(k/2)+(k/2)sin((pi)(x/j)-(pi/2))"
$ws.Range("G3").Value = "(j/2)+j(2*asin(2(y/k)-1)+pi)/2pi"

# ---------------------------------------------------------------------------
# Row 4: Quadratic
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "k(x/j)^2"
$ws.Range("C4").Value = "j*sqrt((y/k))"
$ws.Range("D4").Value = "k-k((x/j)-1)^2"
$ws.Range("E4").Value = "j-j*sqrt(1-(y/k))"
$ws.Range("F4").Value = "2k(x/j)^2 & k-2k((x/j)-1)^2"
$ws.Range("G4").Value = "j*sqrt((y/k)/2) & j-j*sqrt((1-(y/k))/2)"

# ---------------------------------------------------------------------------
# Row 5: Cubic
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "k(x/j)^3"
$ws.Range("C5").Value = "j(y/k)^1/3"
$ws.Range("D5").Value = "k-k(1-(x/j))^3"
$ws.Range("E5").Value = "j-j(1-(y/k))^1/3"
$ws.Range("F5").Value = "4k(x/j)^3 & k-4k((x/j)-1)^3"
$ws.Range("G5").Value = "j((y/k)/4)^1/3 & j-j((1-(y/k))/4)^1/3"

# ---------------------------------------------------------------------------
# Row 6: Quartic
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = "k(x/j)^4"
$ws.Range("C6").Value = "j(y/k)^1/4"
$ws.Range("D6").Value = "k-k((x/j)-1)^4"
$ws.Range("E6").Value = "j-j(1-(y/k))^1/4"
$ws.Range("F6").Value = "8k(x/j)^4 & k-8k((x/j)-1)^4"
$ws.Range("G6").Value = "j((y/k)/8)^1/4 & j-j((1-(y/k))/8)^1/4"

# ---------------------------------------------------------------------------
# Row 7: Quintic
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "k(x/j)^5"
$ws.Range("C7").Value = "j(y/k)^1/5"
$ws.Range("D7").Value = "k-k(1-(x/j))^5"
$ws.Range("E7").Value = "j-j(1-(y/k))^1/5"
$ws.Range("F7").Value = "16k(x/j)^5 & k-16k((x/j)-1)^5"
$ws.Range("G7").Value = "j((y/k)/16)^1/5 & j-j((1-(y/k))/16)^1/5"

# ---------------------------------------------------------------------------
# Row 8: Circular
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "k-k*sqrt(1-(x/j)^2)"
$ws.Range("C8").Value = "j*sqrt(2(y/k)-(y/k)^2)"
$ws.Range("D8").Value = "k*sqrt(2(x/j)-(x/j)^2)"
$ws.Range("E8").Value = "j-j*sqrt(1-(y/k)^2)"
$ws.Range("F8").Value = "(k/2)-(k/2)sqrt(1-4(x/j)^2) & (k/2)-(k/2)sqrt(1-4((x/j)-1)^2)"
$ws.Range("G8").Value = "j*sqrt((y/k)-(y/k)^2) & j-j*sqrt((y/k)-(y/k)^2)"

# ---------------------------------------------------------------------------
# New "Scale x:" / "Scale y:" rows (the new vector replacing the old
# point/origin/destination inputs) - rows 10 & 11, columns A & B.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Scale x:"
$ws.Range("B10").Value = "j"
$ws.Range("A11").Value = "Scale y:"
$ws.Range("B11").Value = "k"

$scaleLabels = $ws.Range("A10:B11")
$scaleLabels.Font.Name = "Consolas"
$scaleLabels.Font.Size = 12
$scaleLabels.Font.Bold = $false
$scaleLabels.NumberFormat = "@"

$ws.Rows.Item(10).RowHeight = 15.75
$ws.Rows.Item(11).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Column widths - re-fit now that each "In/Out/..." column holds differently
# sized text (B..E got narrower, F/G got wider to fit the longest formulas).
# (small constant offset compensates for this engine's internal px rounding
# so the stored character-width lands as close as possible to Excel's
# best-fit output)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 37.59
$ws.Columns.Item(3).ColumnWidth = 34.73
$ws.Columns.Item(4).ColumnWidth = 30.59
$ws.Columns.Item(5).ColumnWidth = 26.45
$ws.Columns.Item(6).ColumnWidth = 85.16
$ws.Columns.Item(7).ColumnWidth = 65.59

# ---------------------------------------------------------------------------
# Selection cursor ends on C9 (blank row under the table) per the saved file.
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()
